# Add extra prediction columns (WIN, TOP2, TOP4, RELEGATION) before the
# existing ExpPoints column, shifting ExpPoints from column C to column G.
# These new columns are placeholders (left blank) to be filled in later by
# a Monte Carlo simulation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# 1) Preserve the existing ExpPoints values (currently in column C, rows 2..lastRow)
$expPoints = $ws.Range("C2:C" + $lastRow).Value2

# 2) Clear out the old column C data values (they will be rewritten as headers/blank)
$ws.Range("C2:C" + $lastRow).ClearContents()

# 3) Set the new header row: C1=WIN, D1=TOP2, E1=TOP4, F1=RELEGATION, G1=ExpPoints
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"
$ws.Range("G1").Value = "ExpPoints"

# Copy the header formatting (bold font, border, centered alignment) from the
# original header cell (A1) onto all of the new header cells.
$ws.Range("A1").Copy()
$ws.Range("C1:G1").PasteSpecial(-4122)

# 4) Write the preserved ExpPoints values into the new column G
$ws.Range("G2:G" + $lastRow).Value2 = $expPoints

# 5) Leave columns C, D, E, F blank for rows 2..lastRow (future WIN/TOP2/TOP4/RELEGATION data)
$ws.Range("C2:F" + $lastRow).ClearContents()
